# "store data in deffrent excel file" - the per-user profile fields
# (name / email / password) move out of this sheet into another workbook;
# only the user id and the face-embedding vector stay here, and the
# column headers are condensed (no more spaces).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: header labels
$ws.Range("A1").Value = "userid"
$ws.Range("B1").Value = "username"
$ws.Range("C1").Value = "useremail"
$ws.Range("D1").Value = "password"
$ws.Range("E1").Value = "userfaceUrl"

# Row 2: sample record - name/email/password no longer live here
$ws.Range("A2").Value = "9zk7nubbtluw8rv5q"
$ws.Range("B2").Value = ""
$ws.Range("C2").Value = ""
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = '{"0":-0.18322962522506714,"1":0.2064070850610733,"2":0.026015130802989006,"3":-0.05736449360847473,"4":-0.022240174934267998,"5":-0.01272329967468977,"6":-0.036502182483673096,"7":-0.08498280495405197,"8":0.16371192038059235,"9":-0.010632830671966076,"10":0.2515646815299988,"11":0.04791298508644104,"12":-0.1866496056318283,"13":-0.07077886909246445,"14":-0.07061983644962311,"15":0.08703090995550156,"16":-0.21615058183670044,"17":-0.16546568274497986,"18":-0.08598027378320694,"19":-0.13444572687149048,"20":0.04641152545809746,"21":0.05617311969399452,"22":-0.03414323553442955,"23":0.035565197467803955,"24":-0.20974335074424744,"25":-0.27684301137924194,"26":-0.02352028526365757,"27":-0.09092475473880768,"28":0.09704017639160156,"29":-0.09487919509410858,"30":0.045347873121500015,"31":-0.03364923223853111,"32":-0.2097923457622528,"33":-0.044465724378824234,"34":0.0443798303604126,"35":0.06176907196640968,"36":-0.001010522129945457,"37":-0.05539001151919365,"38":0.2388012856245041,"39":0.000999712967313826,"40":-0.10115726292133331,"41":0.03114108182489872,"42":0.11585049331188202,"43":0.3301307260990143,"44":0.13646385073661804,"45":0.012258036993443966,"46":0.012249213643372059,"47":-0.06931046396493912,"48":0.0798402950167656,"49":-0.18870915472507477,"50":0.07875624299049377,"51":0.1805676519870758,"52":0.07763119041919708,"53":0.06669747829437256,"54":0.05236174166202545,"55":-0.19541777670383453,"56":0.010519624687731266,"57":0.06189076229929924,"58":-0.22462017834186554,"59":0.05753722041845322,"60":0.044551555067300797,"61":-0.08310028165578842,"62":-0.08822333812713623,"63":-0.0055792140774428844,"64":0.19803865253925323,"65":0.13159707188606262,"66":-0.10665751248598099,"67":-0.201014444231987,"68":0.14996105432510376,"69":-0.14132849872112274,"70":-0.05365066975355148,"71":0.03611328452825546,"72":-0.048990726470947266,"73":-0.10824505984783173,"74":-0.31054240465164185,"75":0.08440668880939484,"76":0.425926148891449,"77":0.12999191880226135,"78":-0.16644300520420074,"79":0.11551357060670853,"80":-0.07487194240093231,"81":-0.06143009290099144,"82":0.06878103315830231,"83":0.053220782428979874,"84":-0.1736316829919815,"85":0.016521165147423744,"86":-0.11116946488618851,"87":0.06857948750257492,"88":0.1636417657136917,"89":0.11347392201423645,"90":-0.0938216969370842,"91":0.1443665325641632,"92":0.026207055896520615,"93":-0.024318447336554527,"94":0.09149383753538132,"95":0.06562624126672745,"96":-0.19301852583885193,"97":-0.03728647530078888,"98":-0.09898153692483902,"99":-0.03784303367137909,"100":0.06756953150033951,"101":-0.039191149175167084,"102":0.03764297813177109,"103":0.16627594828605652,"104":-0.1492624431848526,"105":0.1272362917661667,"106":0.026514194905757904,"107":-0.047301582992076874,"108":-0.02198607847094536,"109":0.09038128703832626,"110":-0.08574368059635162,"111":-0.07558047026395798,"112":0.09492600709199905,"113":-0.24732257425785065,"114":0.21907185018062592,"115":0.11005926132202148,"116":0.03679228201508522,"117":0.17262326180934906,"118":0.15408888459205627,"119":0.037930313497781754,"120":0.015967775136232376,"121":0.07823996245861053,"122":-0.10149620473384857,"123":-0.10199502855539322,"124":-0.03092358447611332,"125":-0.03962467238306999,"126":0.16025269031524658,"127":0.05673786625266075}'
